# Fruta / hortaliza, semanal
# The weekly date-blocks of rows 2-33 (each sharing the same Fecha in column D)
# were reordered. Every column (A:T) for a given source row moves intact to a
# new target row - no individual cell values are altered, only their row
# position. Build the permutation (target row -> source row) and apply it by
# first snapshotting every source row's full A:T values, then writing them
# back out in the new order so reads never see already-overwritten data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (both refer to the ORIGINAL, pre-edit layout)
$mapping = @{
    2  = 9
    3  = 10
    4  = 11
    5  = 12
    6  = 28
    7  = 29
    8  = 30
    9  = 25
    10 = 26
    11 = 27
    12 = 31
    13 = 32
    14 = 33
    15 = 13
    16 = 14
    17 = 15
    18 = 16
    19 = 17
    20 = 18
    21 = 21
    22 = 22
    23 = 23
    24 = 24
    25 = 5
    26 = 6
    27 = 7
    28 = 8
    29 = 2
    30 = 3
    31 = 4
    32 = 19
    33 = 20
}

# Snapshot every source row (2-33), columns A:T, before any writes happen.
$snapshot = @{}
for ($r = 2; $r -le 33; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:T$r").Value2
}

# Write each target row from the snapshot of its mapped source row.
foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $ws.Range("A$targetRow`:T$targetRow").Value2 = $snapshot[$sourceRow]
}
